$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: Caso 6265
$ws.Range("A6").NumberFormat = "@"
$ws.Range("A6").Value = '6265'
$ws.Range("A6").Style = "Normal"
$ws.Range("B6").NumberFormat = "@"
$ws.Range("B6").Value = '8/7/2025'
$ws.Range("B6").Style = "Normal"
$ws.Range("C6").NumberFormat = "@"
$ws.Range("C6").Value = 'BROWN, ALTE. AV. 881'
$ws.Range("C6").Style = "Normal"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = 'ICD30345482'
$ws.Range("E6").Style = "Normal"
$ws.Range("F6").NumberFormat = "@"
$ws.Range("F6").Value = 'Optical Power'
$ws.Range("F6").Style = "Normal"
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = 'Pendiente'
$ws.Range("G6").Style = "Normal"
$ws.Range("H6").NumberFormat = "@"
$ws.Range("H6").Value = 'tendido a baja altura'
$ws.Range("H6").Style = "Normal"
$ws.Range("J6").NumberFormat = "@"
$ws.Range("J6").Value = '{"direccionesNormalizadas": [{"altura": 881, "cod_calle": 2115, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.360551", "y": "-34.632684"}, "direccion": "BROWN, ALTE. AV. 881, CABA", "nombre_calle": "BROWN, ALTE. AV.", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Range("J6").Style = "Normal"
$ws.Range("M6").NumberFormat = "@"
$ws.Range("M6").Value = 'San Telmo'
$ws.Range("M6").Style = "Normal"
$ws.Range("N6").NumberFormat = "@"
$ws.Range("N6").Value = 'Capital Sur'
$ws.Range("N6").Style = "Normal"
$ws.Range("D6").Value = 4
$ws.Range("I6").Value = 1
$ws.Range("K6").Value = -58.360551
$ws.Range("L6").Value = -34.632684

# Row 7: Caso 6899
$ws.Range("A7").NumberFormat = "@"
$ws.Range("A7").Value = '6899'
$ws.Range("A7").Style = "Normal"
$ws.Range("B7").NumberFormat = "@"
$ws.Range("B7").Value = '8/12/2025'
$ws.Range("B7").Style = "Normal"
$ws.Range("C7").NumberFormat = "@"
$ws.Range("C7").Value = 'ERCILLA 6159'
$ws.Range("C7").Style = "Normal"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = 'Pendiente ADM'
$ws.Range("E7").Style = "Normal"
$ws.Range("F7").NumberFormat = "@"
$ws.Range("F7").Value = 'Optical Power'
$ws.Range("F7").Style = "Normal"
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = 'Pendiente'
$ws.Range("G7").Style = "Normal"
$ws.Range("H7").NumberFormat = "@"
$ws.Range("H7").Value = 'Tendido a baja altura y cables cortados'
$ws.Range("H7").Style = "Normal"
$ws.Range("J7").NumberFormat = "@"
$ws.Range("J7").Value = '{"direccionesNormalizadas": [{"altura": 6159, "cod_calle": 5065, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.510283", "y": "-34.647321"}, "direccion": "ERCILLA 6159, CABA", "nombre_calle": "ERCILLA", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Range("J7").Style = "Normal"
$ws.Range("M7").NumberFormat = "@"
$ws.Range("M7").Value = 'Devoto'
$ws.Range("M7").Style = "Normal"
$ws.Range("N7").NumberFormat = "@"
$ws.Range("N7").Value = 'Capital Norte'
$ws.Range("N7").Style = "Normal"
$ws.Range("D7").Value = 9
$ws.Range("I7").Value = 1
$ws.Range("K7").Value = -58.510283
$ws.Range("L7").Value = -34.647321

# Row 8: Caso 6087
$ws.Range("A8").NumberFormat = "@"
$ws.Range("A8").Value = '6087'
$ws.Range("A8").Style = "Normal"
$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = '8/12/2025'
$ws.Range("B8").Style = "Normal"
$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = 'ROSETI 436'
$ws.Range("C8").Style = "Normal"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = 'Pendiente ADM'
$ws.Range("E8").Style = "Normal"
$ws.Range("F8").NumberFormat = "@"
$ws.Range("F8").Value = 'Optical Power'
$ws.Range("F8").Style = "Normal"
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = 'Pendiente'
$ws.Range("G8").Style = "Normal"
$ws.Range("H8").NumberFormat = "@"
$ws.Range("H8").Value = 'Cable en panza'
$ws.Range("H8").Style = "Normal"
$ws.Range("J8").NumberFormat = "@"
$ws.Range("J8").Value = '{"direccionesNormalizadas": [{"altura": 436, "cod_calle": 19088, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.450100", "y": "-34.584687"}, "direccion": "ROSETI 436, CABA", "nombre_calle": "ROSETI", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Range("J8").Style = "Normal"
$ws.Range("M8").NumberFormat = "@"
$ws.Range("M8").Value = 'Colegiales'
$ws.Range("M8").Style = "Normal"
$ws.Range("N8").NumberFormat = "@"
$ws.Range("N8").Value = 'Capital Norte'
$ws.Range("N8").Style = "Normal"
$ws.Range("D8").Value = 15
$ws.Range("I8").Value = 1
$ws.Range("K8").Value = -58.4501
$ws.Range("L8").Value = -34.584687

# Row 9: Caso 6913
$ws.Range("A9").NumberFormat = "@"
$ws.Range("A9").Value = '6913'
$ws.Range("A9").Style = "Normal"
$ws.Range("B9").NumberFormat = "@"
$ws.Range("B9").Value = '8/12/2025'
$ws.Range("B9").Style = "Normal"
$ws.Range("C9").NumberFormat = "@"
$ws.Range("C9").Value = 'RIVADAVIA AV. 6161'
$ws.Range("C9").Style = "Normal"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = 'Pendiente ADM'
$ws.Range("E9").Style = "Normal"
$ws.Range("F9").NumberFormat = "@"
$ws.Range("F9").Value = 'Optical Power'
$ws.Range("F9").Style = "Normal"
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = 'Pendiente'
$ws.Range("G9").Style = "Normal"
$ws.Range("H9").NumberFormat = "@"
$ws.Range("H9").Value = 'Tendido a baja altura'
$ws.Range("H9").Style = "Normal"
$ws.Range("J9").NumberFormat = "@"
$ws.Range("J9").Value = '{"direccionesNormalizadas": [{"altura": 6161, "cod_calle": 19046, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.453305", "y": "-34.625256"}, "direccion": "RIVADAVIA AV. 6161, CABA", "nombre_calle": "RIVADAVIA AV.", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Range("J9").Style = "Normal"
$ws.Range("M9").NumberFormat = "@"
$ws.Range("M9").Value = 'Boedo'
$ws.Range("M9").Style = "Normal"
$ws.Range("N9").NumberFormat = "@"
$ws.Range("N9").Value = 'Capital Sur'
$ws.Range("N9").Style = "Normal"
$ws.Range("D9").Value = 6
$ws.Range("I9").Value = 1
$ws.Range("K9").Value = -58.453305
$ws.Range("L9").Value = -34.625256

# Row 10: Caso 6914
$ws.Range("A10").NumberFormat = "@"
$ws.Range("A10").Value = '6914'
$ws.Range("A10").Style = "Normal"
$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = '8/12/2025'
$ws.Range("B10").Style = "Normal"
$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = 'AVELLANEDA AV. 1240'
$ws.Range("C10").Style = "Normal"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = 'Pendiente ADM'
$ws.Range("E10").Style = "Normal"
$ws.Range("F10").NumberFormat = "@"
$ws.Range("F10").Value = 'Optical Power'
$ws.Range("F10").Style = "Normal"
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = 'Pendiente'
$ws.Range("G10").Style = "Normal"
$ws.Range("H10").NumberFormat = "@"
$ws.Range("H10").Value = 'Tendido a baja altura'
$ws.Range("H10").Style = "Normal"
$ws.Range("J10").NumberFormat = "@"
$ws.Range("J10").Value = '{"direccionesNormalizadas": [{"altura": 1240, "cod_calle": 1141, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.448578", "y": "-34.618227"}, "direccion": "AVELLANEDA AV. 1240, CABA", "nombre_calle": "AVELLANEDA AV.", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Range("J10").Style = "Normal"
$ws.Range("M10").NumberFormat = "@"
$ws.Range("M10").Value = 'Boedo'
$ws.Range("M10").Style = "Normal"
$ws.Range("N10").NumberFormat = "@"
$ws.Range("N10").Value = 'Capital Sur'
$ws.Range("N10").Style = "Normal"
$ws.Range("D10").Value = 6
$ws.Range("I10").Value = 1
$ws.Range("K10").Value = -58.448578
$ws.Range("L10").Value = -34.618227

# Row 11: Caso 6917
$ws.Range("A11").NumberFormat = "@"
$ws.Range("A11").Value = '6917'
$ws.Range("A11").Style = "Normal"
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = '8/12/2025'
$ws.Range("B11").Style = "Normal"
$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = 'BRIN, MINISTRO 1375'
$ws.Range("C11").Style = "Normal"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = 'Pendiente ADM'
$ws.Range("E11").Style = "Normal"
$ws.Range("F11").NumberFormat = "@"
$ws.Range("F11").Value = 'Optical Power'
$ws.Range("F11").Style = "Normal"
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = 'Pendiente'
$ws.Range("G11").Style = "Normal"
$ws.Range("H11").NumberFormat = "@"
$ws.Range("H11").Value = 'tendido a baja altura tiene pendiente tambien una columna en l mismo lugar ver en conjunto'
$ws.Range("H11").Style = "Normal"
$ws.Range("J11").NumberFormat = "@"
$ws.Range("J11").Value = '{"direccionesNormalizadas": [{"altura": 1375, "cod_calle": 2114, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.355342", "y": "-34.635650"}, "direccion": "BRIN, MINISTRO 1375, CABA", "nombre_calle": "BRIN, MINISTRO", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Range("J11").Style = "Normal"
$ws.Range("M11").NumberFormat = "@"
$ws.Range("M11").Value = 'San Telmo'
$ws.Range("M11").Style = "Normal"
$ws.Range("N11").NumberFormat = "@"
$ws.Range("N11").Value = 'Capital Sur'
$ws.Range("N11").Style = "Normal"
$ws.Range("D11").Value = 4
$ws.Range("I11").Value = 1
$ws.Range("K11").Value = -58.355342
$ws.Range("L11").Value = -34.63565

# Row 12: Caso 6933
$ws.Range("A12").NumberFormat = "@"
$ws.Range("A12").Value = '6933'
$ws.Range("A12").Style = "Normal"
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = '8/12/2025'
$ws.Range("B12").Style = "Normal"
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = 'ANDONAEGUI 1894'
$ws.Range("C12").Style = "Normal"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = 'Pendiente ADM'
$ws.Range("E12").Style = "Normal"
$ws.Range("F12").NumberFormat = "@"
$ws.Range("F12").Value = 'Optical Power'
$ws.Range("F12").Style = "Normal"
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = 'Pendiente'
$ws.Range("G12").Style = "Normal"
$ws.Range("H12").NumberFormat = "@"
$ws.Range("H12").Value = 'Ganancia suelta en la columna'
$ws.Range("H12").Style = "Normal"
$ws.Range("J12").NumberFormat = "@"
$ws.Range("J12").Value = '{"direccionesNormalizadas": [{"altura": 1894, "cod_calle": 1082, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.484848", "y": "-34.581325"}, "direccion": "ANDONAEGUI 1894, CABA", "nombre_calle": "ANDONAEGUI", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Range("J12").Style = "Normal"
$ws.Range("M12").NumberFormat = "@"
$ws.Range("M12").Value = 'Paternal'
$ws.Range("M12").Style = "Normal"
$ws.Range("N12").NumberFormat = "@"
$ws.Range("N12").Value = 'Capital Norte'
$ws.Range("N12").Style = "Normal"
$ws.Range("D12").Value = 12
$ws.Range("I12").Value = 1
$ws.Range("K12").Value = -58.484848
$ws.Range("L12").Value = -34.581325

# Row 13: Caso 6937
$ws.Range("A13").NumberFormat = "@"
$ws.Range("A13").Value = '6937'
$ws.Range("A13").Style = "Normal"
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = '8/12/2025'
$ws.Range("B13").Style = "Normal"
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = 'ANCHORENA, TOMAS MANUEL DE, DR. 1566'
$ws.Range("C13").Style = "Normal"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = 'Pendiente ADM'
$ws.Range("E13").Style = "Normal"
$ws.Range("F13").NumberFormat = "@"
$ws.Range("F13").Value = 'Optical Power'
$ws.Range("F13").Style = "Normal"
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = 'Pendiente'
$ws.Range("G13").Style = "Normal"
$ws.Range("H13").NumberFormat = "@"
$ws.Range("H13").Value = 'Cable en panza'
$ws.Range("H13").Style = "Normal"
$ws.Range("J13").NumberFormat = "@"
$ws.Range("J13").Value = '{"direccionesNormalizadas": [{"altura": 1566, "cod_calle": 1078, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.403659", "y": "-34.591997"}, "direccion": "ANCHORENA, TOMAS MANUEL DE, DR. 1566, CABA", "nombre_calle": "ANCHORENA, TOMAS MANUEL DE, DR.", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Range("J13").Style = "Normal"
$ws.Range("M13").NumberFormat = "@"
$ws.Range("M13").Value = 'Recoleta'
$ws.Range("M13").Style = "Normal"
$ws.Range("N13").NumberFormat = "@"
$ws.Range("N13").Value = 'Capital Sur'
$ws.Range("N13").Style = "Normal"
$ws.Range("D13").Value = 2
$ws.Range("I13").Value = 1
$ws.Range("K13").Value = -58.403659
$ws.Range("L13").Value = -34.591997

# Row 14: Caso 6940
$ws.Range("A14").NumberFormat = "@"
$ws.Range("A14").Value = '6940'
$ws.Range("A14").Style = "Normal"
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = '8/12/2025'
$ws.Range("B14").Style = "Normal"
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = 'PERDRIEL 1490'
$ws.Range("C14").Style = "Normal"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = 'Pendiente ADM'
$ws.Range("E14").Style = "Normal"
$ws.Range("F14").NumberFormat = "@"
$ws.Range("F14").Value = 'Optical Power'
$ws.Range("F14").Style = "Normal"
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = 'Pendiente'
$ws.Range("G14").Style = "Normal"
$ws.Range("H14").NumberFormat = "@"
$ws.Range("H14").Value = 'Tendido a baja altura'
$ws.Range("H14").Style = "Normal"
$ws.Range("J14").NumberFormat = "@"
$ws.Range("J14").Value = '{"direccionesNormalizadas": [{"altura": 1490, "cod_calle": 17057, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.385244", "y": "-34.651701"}, "direccion": "PERDRIEL 1490, CABA", "nombre_calle": "PERDRIEL", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Range("J14").Style = "Normal"
$ws.Range("M14").NumberFormat = "@"
$ws.Range("M14").Value = 'San Telmo'
$ws.Range("M14").Style = "Normal"
$ws.Range("N14").NumberFormat = "@"
$ws.Range("N14").Value = 'Capital Sur'
$ws.Range("N14").Style = "Normal"
$ws.Range("D14").Value = 4
$ws.Range("I14").Value = 1
$ws.Range("K14").Value = -58.385244
$ws.Range("L14").Value = -34.651701
